$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 17:54"

# Country name changes (ranking reorder swaps due to updated totals)
$ws.Range("A22").Value = "Italia"
$ws.Range("A23").Value = "Pakistan"
$ws.Range("A62").Value = "Moldavia"
$ws.Range("A63").Value = "Argelia"
$ws.Range("A86").Value = "Grecia"
$ws.Range("A87").Value = "Republica de Macedonia"
$ws.Range("A107").Value = "Jordania"
$ws.Range("A108").Value = "Haiti"
$ws.Range("A109").Value = "Gabon"
$ws.Range("A135").Value = "Reunion"
$ws.Range("A136").Value = "Aruba"
$ws.Range("A137").Value = "Bahamas"
$ws.Range("A160").Value = "Republica de Chipre"
$ws.Range("A161").Value = "Letonia"
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# Updated numeric statistics per country row
$ws.Range("B4").Value = 7325037
$ws.Range("C4").Value = 3694
$ws.Range("D4").Value = 4571236
$ws.Range("E4").Value = 2544299
$ws.Range("G4").Value = 49
$ws.Range("H4").Value = 209502

$ws.Range("B5").Value = 6087454
$ws.Range("C5").Value = 14106
$ws.Range("D5").Value = 5025815
$ws.Range("E5").Value = 965961
$ws.Range("G5").Value = 104
$ws.Range("H5").Value = 95678

$ws.Range("B17").Value = 439013
$ws.Range("C17").Value = 4044
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 42001

$ws.Range("B22").Value = 311364
$ws.Range("C22").Value = 1494
$ws.Range("D22").Value = 225190
$ws.Range("E22").Value = 50323
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = 35851

$ws.Range("B23").Value = 310841
$ws.Range("C23").Value = 566
$ws.Range("D23").Value = 296022
$ws.Range("E23").Value = 8353
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 6466

$ws.Range("B25").Value = 287616
$ws.Range("C25").Value = 1278
$ws.Range("E25").Value = 27277
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 9539

$ws.Range("B29").Value = 154575
$ws.Range("C29").Value = 1450
$ws.Range("D29").Value = 131429
$ws.Range("E29").Value = 13876
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 9270

$ws.Range("B37").Value = 111666
$ws.Range("C37").Value = 280
$ws.Range("D37").Value = 86422
$ws.Range("E37").Value = 23146
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 2098

$ws.Range("B45").Value = 90263
$ws.Range("C45").Value = 171
$ws.Range("D45").Value = 79067
$ws.Range("E45").Value = 7962
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = 3234

$ws.Range("D59").Value = 57393
$ws.Range("E59").Value = 295

$ws.Range("B62").Value = 51194
$ws.Range("C62").Value = 319
$ws.Range("D62").Value = 38217
$ws.Range("E62").Value = 11676
$ws.Range("G62").Value = 14
$ws.Range("H62").Value = 1301

$ws.Range("B63").Value = 51067
$ws.Range("D63").Value = 35860
$ws.Range("E63").Value = 13493
$ws.Range("H63").Value = 1714

$ws.Range("B86").Value = 17707
$ws.Range("C86").Value = 263
$ws.Range("D86").Value = 9989
$ws.Range("E86").Value = 7335
$ws.Range("H86").Value = 383

$ws.Range("B87").Value = 17674
$ws.Range("C87").Value = 45
$ws.Range("D87").Value = 14642
$ws.Range("E87").Value = 2303
$ws.Range("G87").Value = 4
$ws.Range("H87").Value = 729

$ws.Range("B95").Value = 13391
$ws.Range("C95").Value = 132
$ws.Range("D95").Value = 7629
$ws.Range("E95").Value = 5382
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 380

$ws.Range("B101").Value = 10441
$ws.Range("C101").Value = 128
$ws.Range("D101").Value = 6720
$ws.Range("E101").Value = 3558
$ws.Range("G101").Value = 5
$ws.Range("H101").Value = 163

$ws.Range("B103").Value = 9895
$ws.Range("C103").Value = 32
$ws.Range("D103").Value = 9536
$ws.Range("E103").Value = 294

$ws.Range("B107").Value = 9226
$ws.Range("C107").Value = 734
$ws.Range("D107").Value = 4359
$ws.Range("E107").Value = 4816
$ws.Range("G107").Value = 6
$ws.Range("H107").Value = 51

$ws.Range("B108").Value = 8740
$ws.Range("D108").Value = 6688
$ws.Range("E108").Value = 1825
$ws.Range("H108").Value = 227

$ws.Range("B109").Value = 8728
$ws.Range("D109").Value = 7934
$ws.Range("E109").Value = 740
$ws.Range("H109").Value = 54

$ws.Range("B110").Value = 8376
$ws.Range("C110").Value = 19
$ws.Range("D110").Value = 7128
$ws.Range("E110").Value = 1124

$ws.Range("B115").Value = 6170
$ws.Range("C115").Value = 153
$ws.Range("D115").Value = 1741
$ws.Range("E115").Value = 4336
$ws.Range("G115").Value = 4
$ws.Range("H115").Value = 93

$ws.Range("B119").Value = 5483
$ws.Range("C119").Value = 26
$ws.Range("D119").Value = 4787
$ws.Range("E119").Value = 574

$ws.Range("B135").Value = 3882
$ws.Range("C135").Value = 197
$ws.Range("D135").Value = 2819
$ws.Range("E135").Value = 1049
$ws.Range("G135").Value = 3
$ws.Range("H135").Value = 14

$ws.Range("B136").Value = 3844
$ws.Range("D136").Value = 2948
$ws.Range("E136").Value = 871
$ws.Range("H136").Value = 25

$ws.Range("B137").Value = 3838
$ws.Range("D137").Value = 2005
$ws.Range("E137").Value = 1744
$ws.Range("H137").Value = 89

$ws.Range("B144").Value = 3090
$ws.Range("C144").Value = 4
$ws.Range("D144").Value = 2439
$ws.Range("E144").Value = 521

$ws.Range("B160").Value = 1713
$ws.Range("C160").Value = 17
$ws.Range("D160").Value = 1369
$ws.Range("E160").Value = 322
$ws.Range("H160").Value = 22

$ws.Range("B161").Value = 1697
$ws.Range("C161").Value = 21
$ws.Range("D161").Value = 1304
$ws.Range("E161").Value = 357
$ws.Range("H161").Value = 36
